# Fix URL to template.
#
# The cell C5 on Sheet1 held a hyperlinked shared string
# "https://tpl.ottr.xyz/pizza/0.1/NamedPizza.ttl" (with an accompanying
# <hyperlinks> entry / external relationship). The commit corrects the
# template URL and removes the now-stale hyperlink, leaving C5 as plain
# text "http://tpl.ottr.xyz/pizza/0.1/NamedPizza".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the hyperlink that was attached to C5 (rId1 -> ...NamedPizza.ttl).
$ws.Hyperlinks.Delete()

# Correct the template URL text shown/stored in C5.
$ws.Range("C5").Value = "http://tpl.ottr.xyz/pizza/0.1/NamedPizza"

# The saved view had scrolled/selected C5 (topLeftCell C1, selection C5);
# restore the default scroll position and move the selection to C6, as in
# the edited workbook.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C6").Select()
